$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns store plain text (e.g. "71.618.84", "  +3.38%  ")
# rather than numbers. Force text formatting first so Excel does not
# auto-coerce the new values into numbers (and drop things like trailing
# zeros), then restore the default "Normal" style so no formatting diff
# is introduced versus the original workbook.
$valueRange = $ws.Range("D2:E51")
$valueRange.NumberFormat = "@"

$ws.Range("D2").Value = "71.618.84"
$ws.Range("E2").Value = "  +3.38%  "
$ws.Range("D3").Value = "3.687.54"
$ws.Range("E3").Value = "  +8.57%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "589.07"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("D6").Value = "179.61"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").Value = "3.680.94"
$ws.Range("E7").Value = "  +8.60%  "
$ws.Range("D8").Value = "0.622"
$ws.Range("E8").Value = "  +4.90%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("D11").Value = "0.613"
$ws.Range("E11").Value = "  +4.57%  "
$ws.Range("D12").Value = "49.87"
$ws.Range("E12").Value = "  +3.24%  "
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").Value = "4.286.62"
$ws.Range("E14").Value = "  +8.87%  "
$ws.Range("D15").Value = "683.29"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "9.00"
$ws.Range("E16").Value = "  +4.72%  "
$ws.Range("D17").Value = "3.694.80"
$ws.Range("E17").Value = "  +9.11%  "
$ws.Range("D18").Value = "71.749.73"
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("E19").Value = "  +2.21%  "
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("D21").Value = "11.66"
$ws.Range("E21").Value = "  +3.30%  "
$ws.Range("E22").Value = "  +3.46%  "
$ws.Range("D23").Value = "6.26"
$ws.Range("E23").Value = "  +16.86%  "
$ws.Range("D24").Value = "17.83"
$ws.Range("E24").Value = "  +4.55%  "
$ws.Range("D25").Value = "103.90"
$ws.Range("E25").Value = "  +2.67%  "
$ws.Range("E26").Value = "  +3.95%  "
$ws.Range("D27").Value = "2.84"
$ws.Range("E27").Value = "  +5.25%  "
$ws.Range("D28").Value = "10.20"
$ws.Range("E28").Value = "  +4.80%  "
$ws.Range("D29").Value = "35.43"
$ws.Range("E29").Value = "  +5.81%  "
$ws.Range("D30").Value = "9.22"
$ws.Range("E30").Value = "  +5.41%  "
$ws.Range("E31").Value = "  +6.30%  "
$ws.Range("E32").Value = "  +10.59%  "
$ws.Range("D33").Value = "574.60"
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("D34").Value = "11.30"
$ws.Range("E34").Value = "  +2.47%  "
$ws.Range("D35").Value = "0.109"
$ws.Range("E35").Value = "  +3.87%  "
$ws.Range("D36").Value = "59.44"
$ws.Range("E36").Value = "  +2.56%  "
$ws.Range("D37").Value = "3.770.16"
$ws.Range("E37").Value = "  +4.49%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "0.145"
$ws.Range("E39").Value = "  +3.44%  "
$ws.Range("D40").Value = "0.0₃0775"
$ws.Range("E40").Value = "  +3.99%  "
$ws.Range("D41").Value = "35.39"
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "3.47"
$ws.Range("E42").Value = "  +5.06%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "2.79"
$ws.Range("E43").Value = "  +3.25%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0461"
$ws.Range("E44").Value = "  +8.26%  "
$ws.Range("D45").Value = "0.352"
$ws.Range("E45").Value = "  +4.98%  "
$ws.Range("E46").Value = "  +7.71%  "
$ws.Range("D47").Value = "3.37"
$ws.Range("E47").Value = "  -0.76%  "
$ws.Range("E48").Value = "  +4.26%  "
$ws.Range("D49").Value = "1.42"
$ws.Range("E49").Value = "  +2.39%  "
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "134.26"
$ws.Range("E51").Value = "  +2.46%  "

$valueRange.Style = "Normal"
